$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) contains numeric-looking text (e.g. "1.00", "0.0720",
# "63.580.81") that must be preserved verbatim rather than auto-converted to a
# number by Excel. Mark only the specific Price cells being rewritten as Text
# before assigning their new values.
$ws.Range("D2,D3,D5,D6,D8,D9,D10,D11,D12,D13,D15,D16,D17,D18,D20,D21,D22,D23,D24,D25,D27,D28,D29,D30,D31,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D50,D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '63.580.81'
$ws.Cells.Item(2, 5).Value = '  -1.18%  '
$ws.Cells.Item(3, 4).Value = '3.399.24'
$ws.Cells.Item(3, 5).Value = '  -0.37%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$ws.Cells.Item(5, 4).Value = '567.49'
$ws.Cells.Item(5, 5).Value = '  -0.51%  '
$ws.Cells.Item(6, 4).Value = '155.83'
$ws.Cells.Item(6, 5).Value = '  -0.49%  '
$ws.Cells.Item(7, 5).Value = '  +0.06%  '
$ws.Cells.Item(8, 4).Value = '3.402.05'
$ws.Cells.Item(8, 5).Value = '  -0.45%  '
$ws.Cells.Item(9, 4).Value = '0.569'
$ws.Cells.Item(9, 5).Value = '  -7.46%  '
$ws.Cells.Item(10, 4).Value = '7.27'
$ws.Cells.Item(10, 5).Value = '  +1.15%  '
$ws.Cells.Item(11, 4).Value = '0.118'
$ws.Cells.Item(11, 5).Value = '  -2.76%  '
$ws.Cells.Item(12, 4).Value = '0.425'
$ws.Cells.Item(12, 5).Value = '  -3.53%  '
$ws.Cells.Item(13, 4).Value = '3.988.21'
$ws.Cells.Item(13, 5).Value = '  -0.26%  '
$ws.Cells.Item(14, 5).Value = '  -0.24%  '
$ws.Cells.Item(15, 4).Value = '26.93'
$ws.Cells.Item(15, 5).Value = '  -3.24%  '
$ws.Cells.Item(16, 4).Value = '0.0000172'
$ws.Cells.Item(16, 5).Value = '  -8.02%  '
$ws.Cells.Item(17, 4).Value = '63.645.72'
$ws.Cells.Item(17, 5).Value = '  -1.15%  '
$ws.Cells.Item(18, 4).Value = '3.401.81'
$ws.Cells.Item(18, 5).Value = '  -1.45%  '
$ws.Cells.Item(19, 5).Value = '  -4.07%  '
$ws.Cells.Item(20, 4).Value = '13.53'
$ws.Cells.Item(20, 5).Value = '  -3.00%  '
$ws.Cells.Item(21, 4).Value = '380.80'
$ws.Cells.Item(21, 5).Value = '  +1.67%  '
$ws.Cells.Item(22, 4).Value = '7.71'
$ws.Cells.Item(22, 5).Value = '  -3.15%  '
$ws.Cells.Item(23, 4).Value = '1.00'
$ws.Cells.Item(23, 5).Value = '  +0.23%  '
$ws.Cells.Item(24, 4).Value = '70.97'
$ws.Cells.Item(24, 5).Value = '  -1.72%  '
$ws.Cells.Item(25, 4).Value = '0.516'
$ws.Cells.Item(25, 5).Value = '  -6.34%  '
$ws.Cells.Item(26, 5).Value = '  -3.96%  '
$ws.Cells.Item(27, 4).Value = '9.67'
$ws.Cells.Item(27, 5).Value = '  -4.95%  '
$ws.Cells.Item(28, 4).Value = '0.177'
$ws.Cells.Item(28, 5).Value = '  +0.86%  '
$ws.Cells.Item(29, 4).Value = '1.00'
$ws.Cells.Item(29, 5).Value = '  -0.22%  '
$ws.Cells.Item(30, 4).Value = '6.05'
$ws.Cells.Item(30, 5).Value = '  -0.81%  '
$ws.Cells.Item(31, 4).Value = '1.38'
$ws.Cells.Item(31, 5).Value = '  -6.48%  '
$ws.Cells.Item(32, 5).Value = '  -1.74%  '
$ws.Cells.Item(33, 2).Value = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(33, 4).Value = '22.88'
$ws.Cells.Item(33, 5).Value = '  -0.80%  '
$ws.Cells.Item(34, 2).Value = 'Aptos'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(34, 4).Value = '6.93'
$ws.Cells.Item(34, 5).Value = '  -3.66%  '
$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(35, 4).Value = '1.51'
$ws.Cells.Item(35, 5).Value = '  -5.26%  '
$ws.Cells.Item(36, 2).Value = 'Monero'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(36, 4).Value = '160.32'
$ws.Cells.Item(36, 5).Value = '  -0.23%  '
$ws.Cells.Item(37, 2).Value = 'Mantle'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(37, 4).Value = '0.829'
$ws.Cells.Item(37, 5).Value = '  +8.22%  '
$ws.Cells.Item(38, 2).Value = 'Stacks'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(38, 4).Value = '1.82'
$ws.Cells.Item(38, 5).Value = '  -2.74%  '
$ws.Cells.Item(39, 2).Value = 'EnergySwap'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(39, 4).Value = '26.02'
$ws.Cells.Item(39, 5).Value = '  -2.46%  '
$ws.Cells.Item(40, 2).Value = 'Maker'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(40, 4).Value = '2.797.07'
$ws.Cells.Item(40, 5).Value = '  -1.69%  '
$ws.Cells.Item(41, 2).Value = 'Hedera'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(41, 4).Value = '0.0720'
$ws.Cells.Item(41, 5).Value = '  -5.00%  '
$ws.Cells.Item(42, 2).Value = 'OKB'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(42, 4).Value = '42.83'
$ws.Cells.Item(42, 5).Value = '  +0.12%  '
$ws.Cells.Item(43, 2).Value = 'RenderToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(43, 4).Value = '6.37'
$ws.Cells.Item(43, 5).Value = '  -5.46%  '
$ws.Cells.Item(44, 2).Value = 'Filecoin'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(44, 4).Value = '4.37'
$ws.Cells.Item(44, 5).Value = '  -5.19%  '
$ws.Cells.Item(45, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(45, 4).Value = '25.57'
$ws.Cells.Item(45, 5).Value = '  -2.89%  '
$ws.Cells.Item(46, 2).Value = 'VeChain'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(46, 4).Value = '0.0303'
$ws.Cells.Item(46, 5).Value = '  -3.04%  '
$ws.Cells.Item(47, 2).Value = 'Bittensor'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(47, 4).Value = '323.41'
$ws.Cells.Item(47, 5).Value = '  +2.10%  '
$ws.Cells.Item(48, 2).Value = 'dogwifhat'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(48, 4).Value = '2.29'
$ws.Cells.Item(48, 5).Value = '  +7.40%  '
$ws.Cells.Item(49, 2).Value = 'ONDO'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(49, 4).Value = '1.03'
$ws.Cells.Item(49, 5).Value = '  -4.15%  '
$ws.Cells.Item(50, 2).Value = 'Stellar'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(50, 4).Value = '0.103'
$ws.Cells.Item(50, 5).Value = '  -5.45%  '
$ws.Cells.Item(51, 2).Value = 'Cosmos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(51, 4).Value = '6.29'
$ws.Cells.Item(51, 5).Value = '  -4.41%  '
